$d = $word.ActiveDocument

# --- 1. Fill in the three previously-empty paragraphs in the "Finance" section ---
# (paragraphs 52, 54 and 56 of the 11 blank paragraphs that follow the
#  "Finance" Heading1, counting through $d.Paragraphs)

$p1 = $d.Paragraphs.Item(52).Range
$p1.Text = "Considering the market research made we have looked closely into the selling price of the product and the packages. "
$d.Paragraphs.Item(52).Range.Font.Name = "Century Gothic"

$p2 = $d.Paragraphs.Item(54).Range
$p2.Text = "We have been working towards forecasting the sales revenue of the product also financially analysing the business currently and for the near future."
$d.Paragraphs.Item(54).Range.Font.Name = "Century Gothic"

$p3 = $d.Paragraphs.Item(56).Range
$p3.Text = "Let me start by first talking about the selling prices and then how we came about achieving a selling price for the product and the packages. "
$d.Paragraphs.Item(56).Range.Font.Name = "Century Gothic"

# --- 2. Move the "_GoBack" bookmark to the end of the "The sales Revenue is: " paragraph ---
# (Word re-stamps _GoBack at the location of the most recent edit; since the
#  text edits above are earlier in the document, the true last-edit location
#  for this revision sits right after "The sales Revenue is: ")

$target = $d.Content
$target.Find.Execute("The sales Revenue is: ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$target.InsertAfter("IRONGOBACKMARKER")

$marker = $d.Content
$marker.Find.Execute("IRONGOBACKMARKER", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$d.Bookmarks.Add("_GoBack", $marker)

$marker2 = $d.Content
$marker2.Find.Execute("IRONGOBACKMARKER", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$marker2.Text = ""
